$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monitoramento diário")
$ws.Activate()

# Fill in values for row 15 (day 22)
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 33
$ws.Range("G15").Value = 33
$ws.Range("I15").Value = 33

# Update the active selection as recorded in the sheet view
$ws.Range("O10").Select()
